$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F header "time_taken", styled like the other header cells (copy E1's format)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Per-row "time_taken" timestamps for rows 2-43
$timeTaken = @(
    "2021-10-05 10:51:08.425990",
    "2021-10-05 10:51:08.426002",
    "2021-10-05 10:51:08.426006",
    "2021-10-05 10:51:08.426008",
    "2021-10-05 10:51:08.426011",
    "2021-10-05 10:51:08.426014",
    "2021-10-05 10:51:08.426017",
    "2021-10-05 10:51:08.426019",
    "2021-10-05 10:51:08.426022",
    "2021-10-05 10:51:08.426025",
    "2021-10-05 10:51:08.426027",
    "2021-10-05 10:51:08.426030",
    "2021-10-05 10:51:08.426033",
    "2021-10-05 10:51:08.426035",
    "2021-10-05 10:51:08.426038",
    "2021-10-05 10:51:08.426041",
    "2021-10-05 10:51:08.426044",
    "2021-10-05 10:51:08.426046",
    "2021-10-05 10:51:08.426049",
    "2021-10-05 10:51:08.426052",
    "2021-10-05 10:51:08.426054",
    "2021-10-05 10:51:08.426057",
    "2021-10-05 10:51:08.426059",
    "2021-10-05 10:51:08.426062",
    "2021-10-05 10:51:08.426065",
    "2021-10-05 10:51:08.426067",
    "2021-10-05 10:51:08.426070",
    "2021-10-05 10:51:08.426072",
    "2021-10-05 10:51:08.426075",
    "2021-10-05 10:51:08.426077",
    "2021-10-05 10:51:08.426080",
    "2021-10-05 10:51:08.426082",
    "2021-10-05 10:51:08.426085",
    "2021-10-05 10:51:08.426088",
    "2021-10-05 10:51:08.426090",
    "2021-10-05 10:51:08.426093",
    "2021-10-05 10:51:08.426096",
    "2021-10-05 10:51:08.426098",
    "2021-10-05 10:51:08.426101",
    "2021-10-05 10:51:08.426104",
    "2021-10-05 10:51:08.426107",
    "2021-10-05 10:51:08.426109"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
